$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (column F) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6704
$ws1.Range("F4").Value = 424
$ws1.Range("F9").Value = 92
$ws1.Range("F15").Value = 1452
$ws1.Range("F17").Value = 3354
$ws1.Range("F19").Value = 223
$ws1.Range("F21").Value = 1997
$ws1.Range("F22").Value = 110
$ws1.Range("F25").Value = 131

# Sheet "全部类型" (sheet4) - update "想去人数" (column F) counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6704
$ws4.Range("F4").Value = 424
$ws4.Range("F10").Value = 92
$ws4.Range("F16").Value = 1452
$ws4.Range("F18").Value = 3354
$ws4.Range("F20").Value = 223
$ws4.Range("F22").Value = 1997
$ws4.Range("F23").Value = 110
$ws4.Range("F26").Value = 131
